$d = $word.ActiveDocument

# Merge the split runs "[" + "ContractName" + "]" into a single run
# containing the literal text "[ContractName]" (heading paragraph).
$d.Content.Find.Execute("[ContractName]", $true, $false, $false, $false, $false, `
                         $false, 1, $false, "[ContractName]", 2)

# Merge the split runs "[" + "ContractDescription" + "]" into a single run
# containing the literal text "[ContractDescription]" (overview paragraph).
$d.Content.Find.Execute("[ContractDescription]", $true, $false, $false, $false, $false, `
                         $false, 1, $false, "[ContractDescription]", 2)
